$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (43) to the table with Sulfur (rhombic) data.
# Shared-string insertion order matters (matches author's entry order):
# B43="S0", C43="S", A43="Sulfur (rhombic)", then D43, E43.
$ws.Range("B43").Value = "S0"
$ws.Range("C43").Value = "S"
$ws.Range("A43").Value = "Sulfur (rhombic)"
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = "Kleerebezem2010"

# Copy the center-aligned style from neighboring cells (B42:E42) to B43:E43
$ws.Range("B42:E42").Copy()
$ws.Range("B43:E43").PasteSpecial(-4122) # xlPasteFormats

# Update the selection to match the final state
$ws.Range("G38").Select()
